$d = $word.ActiveDocument

$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -eq "THE SMARTCARD INTRODUCTION`r") {
        $startPara = $p
    }
    if ($t -like "We*ve prepared an short video*") {
        $endPara = $p
    }
}

if ($startPara -eq $null -or $endPara -eq $null) {
    Write-Host "Could not locate target paragraphs"
} else {
    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)

    $pkgXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading2"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="150" w:afterAutospacing="0" w:line="264" w:lineRule="atLeast"/><w:jc w:val="center"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="exo" w:hAnsi="exo" w:cs="open sans"/><w:caps/><w:color w:val="343434"/><w:spacing w:val="15"/><w:sz w:val="42"/><w:szCs w:val="42"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="exo" w:hAnsi="exo" w:cs="open sans"/><w:caps/><w:color w:val="343434"/><w:spacing w:val="15"/><w:sz w:val="42"/><w:szCs w:val="42"/></w:rPr><w:t xml:space="preserve">THE SMARTCARD INTRODUCTION</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading2"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="150" w:afterAutospacing="0" w:line="264" w:lineRule="atLeast"/><w:jc w:val="center"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="exo" w:hAnsi="exo"/><w:caps/><w:color w:val="343434"/><w:spacing w:val="15"/><w:sz w:val="42"/><w:szCs w:val="42"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="exo" w:hAnsi="exo"/><w:caps/><w:color w:val="343434"/><w:spacing w:val="15"/><w:sz w:val="42"/><w:szCs w:val="42"/></w:rPr><w:t xml:space="preserve">SMARTCASH FOR BUSINESS</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading2"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="150" w:afterAutospacing="0" w:line="264" w:lineRule="atLeast"/><w:jc w:val="center"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="exo" w:hAnsi="exo" w:cs="open sans"/><w:caps/><w:color w:val="343434"/><w:spacing w:val="15"/><w:sz w:val="42"/><w:szCs w:val="42"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Heading1"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="210" w:afterAutospacing="0" w:line="264" w:lineRule="atLeast"/><w:jc w:val="center"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="exo" w:hAnsi="exo"/><w:caps/><w:color w:val="343434"/><w:spacing w:val="15"/><w:sz w:val="51"/><w:szCs w:val="51"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="exo" w:hAnsi="exo"/><w:caps/><w:color w:val="343434"/><w:spacing w:val="15"/><w:sz w:val="51"/><w:szCs w:val="51"/></w:rPr><w:t xml:space="preserve">SMARTCARD TUTORIAL VIDEOS</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading1"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="210" w:afterAutospacing="0" w:line="264" w:lineRule="atLeast"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="exo" w:hAnsi="exo"/><w:caps/><w:color w:val="343434"/><w:spacing w:val="15"/><w:sz w:val="51"/><w:szCs w:val="51"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="exo" w:hAnsi="exo"/><w:caps/><w:color w:val="343434"/><w:spacing w:val="15"/><w:sz w:val="51"/><w:szCs w:val="51"/></w:rPr><w:t xml:space="preserve">DOWNLOAD SMARTPAY APP</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NormalWeb"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="204" w:beforeAutospacing="0" w:after="204" w:afterAutospacing="0"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="open sans" w:hAnsi="open sans" w:cs="open sans"/><w:color w:val="252525"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="open sans" w:hAnsi="open sans" w:cs="open sans"/><w:color w:val="252525"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve">Accept SmartCash as a payment option in your business with zero fees using a simple SmartPay app. The SmartPay app is available for use anywhere in the world, all it requires is access to the internet.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading2"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="150" w:afterAutospacing="0" w:line="264" w:lineRule="atLeast"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="exo" w:hAnsi="exo"/><w:caps/><w:color w:val="343434"/><w:spacing w:val="15"/><w:sz w:val="42"/><w:szCs w:val="42"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="exo" w:hAnsi="exo"/><w:caps/><w:color w:val="343434"/><w:spacing w:val="15"/><w:sz w:val="42"/><w:szCs w:val="42"/></w:rPr><w:t xml:space="preserve">SMARTCARD SHOP</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NormalWeb"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="204" w:beforeAutospacing="0" w:after="204" w:afterAutospacing="0"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="open sans" w:hAnsi="open sans" w:cs="open sans"/><w:color w:val="252525"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="open sans" w:hAnsi="open sans" w:cs="open sans"/><w:color w:val="252525"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve">An online store where you can purchase physical SmartCards.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading2"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="150" w:afterAutospacing="0" w:line="264" w:lineRule="atLeast"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="exo" w:hAnsi="exo"/><w:caps/><w:color w:val="343434"/><w:spacing w:val="15"/><w:sz w:val="42"/><w:szCs w:val="42"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="exo" w:hAnsi="exo"/><w:caps/><w:color w:val="343434"/><w:spacing w:val="15"/><w:sz w:val="42"/><w:szCs w:val="42"/></w:rPr><w:t xml:space="preserve">READY TO GET STARTED?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading2"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="150" w:afterAutospacing="0" w:line="264" w:lineRule="atLeast"/><w:jc w:val="center"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="exo" w:hAnsi="exo" w:cs="open sans"/><w:caps/><w:color w:val="343434"/><w:spacing w:val="15"/><w:sz w:val="42"/><w:szCs w:val="42"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:line="396" w:lineRule="atLeast"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="inherit" w:hAnsi="inherit" w:cs="open sans"/><w:color w:val="252525"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

    $r.InsertXML($pkgXml, $null) | Out-Null
    Write-Host "Replacement complete"
}
